$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 44

# Copy formatting from the last existing data row so the new row matches styling
$ws.Range("A43:D43").Copy($ws.Range("A44:D44"))

$ws.Cells.Item($row, 1).Value = "Z03_B04_P01"
$ws.Cells.Item($row, 2).Value = "Z03_B04"
$ws.Cells.Item($row, 3).Value = "Unterschiede in der mittleren Lebenserwartung zwischen den Lebensräumen mit hoher bzw. niedriger sozioökonomischer Deprivation reduzieren"
$ws.Cells.Item($row, 4).Value = "X"
